$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark from the near-end empty paragraph to
#        the start of the "Peer-reviewed research journals..." bullet. ---

# Remove the existing _GoBack bookmark (wherever Word last left it).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-create it, collapsed, right before "Peer-reviewed research journals..."
$allText = $d.Content.Text
$peerIdx = $allText.IndexOf("Peer-reviewed research journals")
$bookmarkRange = $d.Range($peerIdx, $peerIdx)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 2. Change "backed" -> "proved" in the writing-style sentence, which
#        splits that sentence's run into three runs (same formatting). ---

$allText = $d.Content.Text
$wordIdx = $allText.IndexOf("backed")
$targetRange = $d.Range($wordIdx, $wordIdx + 6)
$targetRange.Text = "proved"

# Force the run split around the replaced word by nudging then reverting a
# character attribute on that exact range (mirrors how Word splits runs
# when formatting is touched, even though the net formatting is unchanged).
$splitRange = $d.Range($wordIdx, $wordIdx + 6)
$splitRange.Bold = 1
$splitRange.Bold = 0
